$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: two new time-log entries (rows 107-108) ---

$ws1.Range("A107").Value = 41950
$ws1.Range("B107").Value = 0.7270833333333333
$ws1.Range("C107").Value = 0.87222222222222223
$ws1.Range("D107").Value = 30
$ws1.Range("E107").Formula = "=IF(AND(NOT(ISBLANK(B107)),NOT(ISBLANK(C107))), (C107-B107) * 24 - D107/60, """")"
$ws1.Range("F107").Value = "Coding"

$ws1.Range("A108").Value = 41951
$ws1.Range("B108").Value = 0.59722222222222221
$ws1.Range("C108").Value = 0.65347222222222223
$ws1.Range("D108").Value = 15
$ws1.Range("E108").Formula = "=IF(AND(NOT(ISBLANK(B108)),NOT(ISBLANK(C108))), (C108-B108) * 24 - D108/60, """")"
$ws1.Range("F108").Value = "Coding"

# Rows 109-120 stay blank in A-D/F, but still carry the (empty-string) delta
# formula like every other still-unused row in the log above row 121.
for ($r = 109; $r -le 120; $r++) {
    $ws1.Range("E$r").Formula = "=IF(AND(NOT(ISBLANK(B$r)),NOT(ISBLANK(C$r))), (C$r-B$r) * 24 - D$r/60, """")"
}

# New trailing row with a "# days" custom format cell, left blank.
$ws1.Range("E123").NumberFormat = '#" days"'

# --- Sheet2: append a "Total:" row under the category breakdown ---

$ws2.Range("A6").Value = "Total:"
$ws2.Range("B6").Formula = "=SUM(B2:B5)"
$ws2.Range("B6").NumberFormat = $ws2.Range("B2").NumberFormat

$ws2.Columns.Item(2).ColumnWidth = 5.7

$r1 = $ws2.Range("A7")
$r2 = $ws2.Range("C9")
$u = $excel.Union($r1, $r2)
$u.Select()

# Re-activate Sheet1 (it is the tab that should stay selected) and restore
# its own selection last, so the final view state sticks to this sheet.
$ws1.Range("D116").Select()
